# Insert a new weekly price record as row 94, pushing the existing
# rows 94:158 down to 95:159 (the last row's data is preserved, now at
# row 159). This matches the commit "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 94:158 down to 95:159, carrying formatting along (mirrors
# Excel's Insert Row UI action).
$ws.Rows(94).Insert()

# Populate the newly inserted row 94 with the new weekly observation.
$ws.Range("A94").Value = 11
$ws.Range("B94").Value = "Vega Monumental Concepción"
$ws.Range("C94").Value = "Bíobío"
$ws.Range("D94").Value = 44827
$ws.Range("E94").Value = 8
$ws.Range("F94").Value = 100112043
$ws.Range("G94").Value = "Pepino ensalada"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 100
$ws.Range("K94").Value = 17000
$ws.Range("L94").Value = 18000
$ws.Range("M94").Value = 17500
$ws.Range("N94").Value = '$/caja 60 unidades'
$ws.Range("O94").Value = "Región de Arica y Parinacota"
$ws.Range("P94").Value = 292
$ws.Range("Q94").Value = 60
$ws.Range("R94").Value = "Hortaliza"
